{"js": "const replacements = [\n  [\"28\u00f77=4, 0\", \"84\u00f75=16, 4\"],\n  [\"34\u00f79=3, 7\", \"90\u00f76=15, 0\"],\n  [\"96\u00f79=10, 6\", \"76\u00f76=12, 4\"],\n  [\"59\u00f76=9, 5\", \"28\u00f79=3, 1\"],\n  [\"96\u00f78=12, 0\", \"65\u00f74=16, 1\"],\n  [\"42\u00f72=21, 0\", \"29\u00f75=5, 4\"],\n  [\"11\u00f75=2, 1\", \"77\u00f77=11, 0\"],\n  [\"12\u00f77=1, 5\", \"27\u00f79=3, 0\"],\n  [\"75\u00f79=8, 3\", \"45\u00f79=5, 0\"],\n  [\"43\u00f72=21, 1\", \"14\u00f73=4, 2\"],\n  [\"30\u00f72=15, 0\", \"77\u00f76=12, 5\"],\n  [\"44\u00f76=7, 2\", \"17\u00f73=5, 2\"],\n  [\"45\u00f76=7, 3\", \"93\u00f79=10, 3\"],\n  [\"25\u00f76=4, 1\", \"89\u00f72=44, 1\"],\n  [\"97\u00f76=16, 1\", \"18\u00f79=2, 0\"],\n  [\"90\u00f73=30, 0\", \"57\u00f79=6, 3\"],\n  [\"85\u00f78=10, 5\", \"58\u00f76=9, 4\"],\n  [\"94\u00f76=15, 4\", \"81\u00f79=9, 0\"],\n  [\"49\u00f74=12, 1\", \"51\u00f78=6, 3\"],\n  [\"12\u00f79=1, 3\", \"38\u00f73=12, 2\"],\n  [\"62\u00f77=8, 6\", \"39\u00f77=5, 4\"],\n  [\"88\u00f76=14, 4\", \"49\u00f72=24, 1\"],\n  [\"81\u00f74=20, 1\", \"96\u00f74=24, 0\"],\n  [\"34\u00f77=4, 6\", \"76\u00f74=19, 0\"],\n  [\"50\u00f76=8, 2\", \"80\u00f76=13, 2\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('28\u00f77=4, 0', '84\u00f75=16, 4')\n    ,@('34\u00f79=3, 7', '90\u00f76=15, 0')\n    ,@('96\u00f79=10, 6', '76\u00f76=12, 4')\n    ,@('59\u00f76=9, 5', '28\u00f79=3, 1')\n    ,@('96\u00f78=12, 0', '65\u00f74=16, 1')\n    ,@('42\u00f72=21, 0', '29\u00f75=5, 4')\n    ,@('11\u00f75=2, 1', '77\u00f77=11, 0')\n    ,@('12\u00f77=1, 5', '27\u00f79=3, 0')\n    ,@('75\u00f79=8, 3', '45\u00f79=5, 0')\n    ,@('43\u00f72=21, 1', '14\u00f73=4, 2')\n    ,@('30\u00f72=15, 0', '77\u00f76=12, 5')\n    ,@('44\u00f76=7, 2', '17\u00f73=5, 2')\n    ,@('45\u00f76=7, 3', '93\u00f79=10, 3')\n    ,@('25\u00f76=4, 1', '89\u00f72=44, 1')\n    ,@('97\u00f76=16, 1', '18\u00f79=2, 0')\n    ,@('90\u00f73=30, 0', '57\u00f79=6, 3')\n    ,@('85\u00f78=10, 5', '58\u00f76=9, 4')\n    ,@('94\u00f76=15, 4', '81\u00f79=9, 0')\n    ,@('49\u00f74=12, 1', '51\u00f78=6, 3')\n    ,@('12\u00f79=1, 3', '38\u00f73=12, 2')\n    ,@('62\u00f77=8, 6', '39\u00f77=5, 4')\n    ,@('88\u00f76=14, 4', '49\u00f72=24, 1')\n    ,@('81\u00f74=20, 1', '96\u00f74=24, 0')\n    ,@('34\u00f77=4, 6', '76\u00f74=19, 0')\n    ,@('50\u00f76=8, 2', '80\u00f76=13, 2')\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"No match found for: $findText\"\n    }\n}"}
